$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update swapped match rows (home/away pairs whose full data rows were exchanged) ---
$ws.Cells.Item(3, 2).Value = 6776469
$ws.Cells.Item(3, 5).Value = 'Magdeburg II'
$ws.Cells.Item(3, 6).Value = 'SG RotWeiss Thalheim'
$ws.Cells.Item(3, 10).Value = 1.05
$ws.Cells.Item(3, 11).Value = 13
$ws.Cells.Item(3, 12).Value = 19
$ws.Cells.Item(3, 13).Value = 1.05
$ws.Cells.Item(3, 14).Value = 13
$ws.Cells.Item(3, 15).Value = 19
$ws.Cells.Item(3, 16).Value = -3
$ws.Cells.Item(3, 17).Value = 1.9
$ws.Cells.Item(3, 18).Value = 1.9
$ws.Cells.Item(3, 19).Value = 3.5
$ws.Cells.Item(3, 20).Value = 1.775
$ws.Cells.Item(3, 21).Value = 1.925
$ws.Cells.Item(3, 23).Value = 12
$ws.Cells.Item(3, 25).Value = -1
$ws.Cells.Item(3, 26).Value = 0.8999999999999999
$ws.Cells.Item(3, 28).Value = 0.925
$ws.Cells.Item(4, 2).Value = 6776470
$ws.Cells.Item(4, 5).Value = 'SV Dessau 05'
$ws.Cells.Item(4, 6).Value = '1 FC BitterfeldWolfen'
$ws.Cells.Item(4, 10).Value = 2.2
$ws.Cells.Item(4, 11).Value = 3.25
$ws.Cells.Item(4, 12).Value = 2.8
$ws.Cells.Item(4, 13).Value = 2.15
$ws.Cells.Item(4, 14).Value = 3.5
$ws.Cells.Item(4, 15).Value = 2.7
$ws.Cells.Item(4, 16).Value = -0.25
$ws.Cells.Item(4, 17).Value = 1.975
$ws.Cells.Item(4, 18).Value = 1.825
$ws.Cells.Item(4, 19).Value = 4
$ws.Cells.Item(4, 20).Value = 1.925
$ws.Cells.Item(4, 21).Value = 1.875
$ws.Cells.Item(4, 23).Value = 2.5
$ws.Cells.Item(4, 25).Value = -0.5
$ws.Cells.Item(4, 26).Value = 0.4125
$ws.Cells.Item(4, 28).Value = 0.875
$ws.Cells.Item(16, 2).Value = 7138608
$ws.Cells.Item(16, 5).Value = 'SV UnterFlockenbach'
$ws.Cells.Item(16, 6).Value = 'SC Dortelweil'
$ws.Cells.Item(16, 7).Value = 1
$ws.Cells.Item(16, 8).Value = 1
$ws.Cells.Item(16, 9).Value = 'D'
$ws.Cells.Item(16, 10).Value = 1.083
$ws.Cells.Item(16, 11).Value = 9
$ws.Cells.Item(16, 12).Value = 16
$ws.Cells.Item(16, 13).Value = 1.125
$ws.Cells.Item(16, 14).Value = 7.5
$ws.Cells.Item(16, 15).Value = 13
$ws.Cells.Item(16, 16).Value = -2.5
$ws.Cells.Item(16, 17).Value = 1.775
$ws.Cells.Item(16, 18).Value = 1.925
$ws.Cells.Item(16, 19).Value = 4.25
$ws.Cells.Item(16, 20).Value = 1.975
$ws.Cells.Item(16, 21).Value = 1.825
$ws.Cells.Item(16, 22).Value = -1
$ws.Cells.Item(16, 23).Value = 6.5
$ws.Cells.Item(16, 25).Value = -1
$ws.Cells.Item(16, 26).Value = 0.925
$ws.Cells.Item(16, 27).Value = -1
$ws.Cells.Item(16, 28).Value = 0.825
$ws.Cells.Item(17, 2).Value = 7138607
$ws.Cells.Item(17, 5).Value = 'Rot Weiss Walldorf II'
$ws.Cells.Item(17, 6).Value = 'Turnerschaft OberRoden'
$ws.Cells.Item(17, 7).Value = 3
$ws.Cells.Item(17, 8).Value = 2
$ws.Cells.Item(17, 9).Value = 'H'
$ws.Cells.Item(17, 10).Value = 2.25
$ws.Cells.Item(17, 11).Value = 3.75
$ws.Cells.Item(17, 12).Value = 2.5
$ws.Cells.Item(17, 13).Value = 2.25
$ws.Cells.Item(17, 14).Value = 3.8
$ws.Cells.Item(17, 15).Value = 2.45
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(17, 17).Value = 1.8
$ws.Cells.Item(17, 18).Value = 2
$ws.Cells.Item(17, 19).Value = 3.75
$ws.Cells.Item(17, 20).Value = 1.95
$ws.Cells.Item(17, 21).Value = 1.85
$ws.Cells.Item(17, 22).Value = 1.25
$ws.Cells.Item(17, 23).Value = -1
$ws.Cells.Item(17, 25).Value = 0.8
$ws.Cells.Item(17, 26).Value = -1
$ws.Cells.Item(17, 27).Value = 0.95
$ws.Cells.Item(17, 28).Value = -1
$ws.Cells.Item(46, 2).Value = 7248441
$ws.Cells.Item(46, 5).Value = '1 FC Lok Stendal'
$ws.Cells.Item(46, 6).Value = 'SV BlauWeiss Zorbau'
$ws.Cells.Item(46, 7).Value = 3
$ws.Cells.Item(46, 8).Value = 1
$ws.Cells.Item(46, 9).Value = 'H'
$ws.Cells.Item(46, 10).Value = 3.25
$ws.Cells.Item(46, 11).Value = 3.8
$ws.Cells.Item(46, 12).Value = 1.833
$ws.Cells.Item(46, 13).Value = 2.7
$ws.Cells.Item(46, 14).Value = 3.75
$ws.Cells.Item(46, 15).Value = 2.1
$ws.Cells.Item(46, 17).Value = 1.85
$ws.Cells.Item(46, 18).Value = 1.95
$ws.Cells.Item(46, 19).Value = 3
$ws.Cells.Item(46, 22).Value = 1.7
$ws.Cells.Item(46, 24).Value = -1
$ws.Cells.Item(46, 25).Value = 0.8500000000000001
$ws.Cells.Item(46, 26).Value = -1
$ws.Cells.Item(47, 2).Value = 7248791
$ws.Cells.Item(47, 5).Value = 'Waldhof Mannheim II'
$ws.Cells.Item(47, 6).Value = 'FC Zuzenhausen'
$ws.Cells.Item(47, 7).Value = 1
$ws.Cells.Item(47, 8).Value = 3
$ws.Cells.Item(47, 9).Value = 'A'
$ws.Cells.Item(47, 10).Value = 2.75
$ws.Cells.Item(47, 11).Value = 4
$ws.Cells.Item(47, 12).Value = 2
$ws.Cells.Item(47, 13).Value = 2.75
$ws.Cells.Item(47, 14).Value = 4
$ws.Cells.Item(47, 15).Value = 2
$ws.Cells.Item(47, 17).Value = 1.975
$ws.Cells.Item(47, 18).Value = 1.825
$ws.Cells.Item(47, 19).Value = 3.25
$ws.Cells.Item(47, 22).Value = -1
$ws.Cells.Item(47, 24).Value = 1
$ws.Cells.Item(47, 25).Value = -1
$ws.Cells.Item(47, 26).Value = 0.825
$ws.Cells.Item(69, 2).Value = 7423699
$ws.Cells.Item(69, 5).Value = 'SG 2000 MulheimKarlich'
$ws.Cells.Item(69, 6).Value = 'Ahrweiler BC'
$ws.Cells.Item(69, 7).Value = 2
$ws.Cells.Item(69, 8).Value = 2
$ws.Cells.Item(69, 9).Value = 'D'
$ws.Cells.Item(69, 10).Value = 2.2
$ws.Cells.Item(69, 11).Value = 5
$ws.Cells.Item(69, 12).Value = 2.2
$ws.Cells.Item(69, 13).Value = 2.2
$ws.Cells.Item(69, 14).Value = 4.75
$ws.Cells.Item(69, 15).Value = 2.2
$ws.Cells.Item(69, 16).Value = 0
$ws.Cells.Item(69, 17).Value = 1.9
$ws.Cells.Item(69, 18).Value = 1.9
$ws.Cells.Item(69, 19).Value = 4.25
$ws.Cells.Item(69, 20).Value = 1.775
$ws.Cells.Item(69, 21).Value = 2.025
$ws.Cells.Item(69, 22).Value = -1
$ws.Cells.Item(69, 23).Value = 3.75
$ws.Cells.Item(69, 25).Value = 0
$ws.Cells.Item(69, 26).Value = 0
$ws.Cells.Item(69, 27).Value = -0.5
$ws.Cells.Item(69, 28).Value = 0.5125
$ws.Cells.Item(70, 2).Value = 7423700
$ws.Cells.Item(70, 5).Value = 'TuS Hornau'
$ws.Cells.Item(70, 6).Value = 'FC Burgsolms'
$ws.Cells.Item(70, 7).Value = 3
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 'H'
$ws.Cells.Item(70, 10).Value = 1.727
$ws.Cells.Item(70, 11).Value = 4.5
$ws.Cells.Item(70, 12).Value = 3.2
$ws.Cells.Item(70, 13).Value = 1.727
$ws.Cells.Item(70, 14).Value = 4.5
$ws.Cells.Item(70, 15).Value = 3.2
$ws.Cells.Item(70, 16).Value = -0.5
$ws.Cells.Item(70, 17).Value = 1.775
$ws.Cells.Item(70, 18).Value = 2.025
$ws.Cells.Item(70, 19).Value = 3.5
$ws.Cells.Item(70, 20).Value = 1.85
$ws.Cells.Item(70, 21).Value = 1.95
$ws.Cells.Item(70, 22).Value = 0.7270000000000001
$ws.Cells.Item(70, 23).Value = -1
$ws.Cells.Item(70, 25).Value = 0.7749999999999999
$ws.Cells.Item(70, 26).Value = -1
$ws.Cells.Item(70, 27).Value = -1
$ws.Cells.Item(70, 28).Value = 0.95
$ws.Cells.Item(86, 2).Value = 7511976
$ws.Cells.Item(86, 5).Value = 'DJK Bad Homburg'
$ws.Cells.Item(86, 6).Value = 'SG Bornheim 1945 GrunWeiss'
$ws.Cells.Item(86, 7).Value = 4
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 10).Value = 2
$ws.Cells.Item(86, 11).Value = 3.75
$ws.Cells.Item(86, 12).Value = 2.9
$ws.Cells.Item(86, 13).Value = 1.8
$ws.Cells.Item(86, 15).Value = 3.3
$ws.Cells.Item(86, 16).Value = -0.5
$ws.Cells.Item(86, 17).Value = 1.85
$ws.Cells.Item(86, 18).Value = 1.95
$ws.Cells.Item(86, 19).Value = 3.5
$ws.Cells.Item(86, 20).Value = 1.975
$ws.Cells.Item(86, 21).Value = 1.825
$ws.Cells.Item(86, 22).Value = 0.8
$ws.Cells.Item(86, 25).Value = 0.8500000000000001
$ws.Cells.Item(86, 27).Value = 0.9750000000000001
$ws.Cells.Item(86, 28).Value = -1
$ws.Cells.Item(87, 2).Value = 7511958
$ws.Cells.Item(87, 5).Value = 'SpVgg EGC Wirges'
$ws.Cells.Item(87, 6).Value = 'SG 2000 MulheimKarlich'
$ws.Cells.Item(87, 7).Value = 2
$ws.Cells.Item(87, 8).Value = 1
$ws.Cells.Item(87, 10).Value = 4.333
$ws.Cells.Item(87, 11).Value = 4
$ws.Cells.Item(87, 12).Value = 1.571
$ws.Cells.Item(87, 13).Value = 4.2
$ws.Cells.Item(87, 15).Value = 1.571
$ws.Cells.Item(87, 16).Value = 1
$ws.Cells.Item(87, 17).Value = 1.875
$ws.Cells.Item(87, 18).Value = 1.925
$ws.Cells.Item(87, 19).Value = 3.75
$ws.Cells.Item(87, 20).Value = 1.925
$ws.Cells.Item(87, 21).Value = 1.875
$ws.Cells.Item(87, 22).Value = 3.2
$ws.Cells.Item(87, 25).Value = 0.875
$ws.Cells.Item(87, 27).Value = -1
$ws.Cells.Item(87, 28).Value = 0.875
$ws.Cells.Item(125, 2).Value = 8039381
$ws.Cells.Item(125, 5).Value = 'SG Union Klosterfelde'
$ws.Cells.Item(125, 6).Value = 'SV 1908 GW Ahrensfelde'
$ws.Cells.Item(125, 7).Value = 1
$ws.Cells.Item(125, 8).Value = 3
$ws.Cells.Item(125, 10).Value = 3.25
$ws.Cells.Item(125, 11).Value = 3.8
$ws.Cells.Item(125, 12).Value = 1.833
$ws.Cells.Item(125, 13).Value = 3.25
$ws.Cells.Item(125, 14).Value = 3.8
$ws.Cells.Item(125, 15).Value = 1.833
$ws.Cells.Item(125, 16).Value = 0.5
$ws.Cells.Item(125, 17).Value = 1.925
$ws.Cells.Item(125, 18).Value = 1.875
$ws.Cells.Item(125, 19).Value = 3
$ws.Cells.Item(125, 20).Value = 1.825
$ws.Cells.Item(125, 21).Value = 1.975
$ws.Cells.Item(125, 24).Value = 0.833
$ws.Cells.Item(125, 26).Value = 0.875
$ws.Cells.Item(125, 27).Value = 0.825
$ws.Cells.Item(125, 28).Value = -1
$ws.Cells.Item(126, 2).Value = 8039382
$ws.Cells.Item(126, 5).Value = 'FC Burgsolms'
$ws.Cells.Item(126, 6).Value = 'TSV Steinbach II'
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 4
$ws.Cells.Item(126, 10).Value = 15
$ws.Cells.Item(126, 11).Value = 9
$ws.Cells.Item(126, 12).Value = 1.111
$ws.Cells.Item(126, 13).Value = 15
$ws.Cells.Item(126, 14).Value = 9
$ws.Cells.Item(126, 15).Value = 1.111
$ws.Cells.Item(126, 16).Value = 2.75
$ws.Cells.Item(126, 17).Value = 1.9
$ws.Cells.Item(126, 18).Value = 1.9
$ws.Cells.Item(126, 19).Value = 4
$ws.Cells.Item(126, 20).Value = 1.9
$ws.Cells.Item(126, 21).Value = 1.9
$ws.Cells.Item(126, 24).Value = 0.111
$ws.Cells.Item(126, 26).Value = 0.8999999999999999
$ws.Cells.Item(126, 27).Value = 0
$ws.Cells.Item(126, 28).Value = 0
$ws.Cells.Item(143, 2).Value = 8121110
$ws.Cells.Item(143, 5).Value = 'FC Astoria Walldorf II'
$ws.Cells.Item(143, 6).Value = 'SV Spielberg'
$ws.Cells.Item(143, 7).Value = 1
$ws.Cells.Item(143, 9).Value = 'D'
$ws.Cells.Item(143, 10).Value = 2
$ws.Cells.Item(143, 11).Value = 3.6
$ws.Cells.Item(143, 12).Value = 3
$ws.Cells.Item(143, 13).Value = 2
$ws.Cells.Item(143, 14).Value = 3.6
$ws.Cells.Item(143, 15).Value = 3
$ws.Cells.Item(143, 16).Value = -0.25
$ws.Cells.Item(143, 17).Value = 1.825
$ws.Cells.Item(143, 18).Value = 1.975
$ws.Cells.Item(143, 19).Value = 3.25
$ws.Cells.Item(143, 20).Value = 1.95
$ws.Cells.Item(143, 21).Value = 1.85
$ws.Cells.Item(143, 22).Value = -1
$ws.Cells.Item(143, 23).Value = 2.6
$ws.Cells.Item(143, 25).Value = -0.5
$ws.Cells.Item(143, 26).Value = 0.4875
$ws.Cells.Item(143, 28).Value = 0.8500000000000001
$ws.Cells.Item(144, 2).Value = 8121117
$ws.Cells.Item(144, 5).Value = 'RotWeiss Frankfurt'
$ws.Cells.Item(144, 6).Value = 'FCA 04 Darmstadt'
$ws.Cells.Item(144, 7).Value = 2
$ws.Cells.Item(144, 9).Value = 'H'
$ws.Cells.Item(144, 10).Value = 2.7
$ws.Cells.Item(144, 11).Value = 3.75
$ws.Cells.Item(144, 12).Value = 2.1
$ws.Cells.Item(144, 13).Value = 2.75
$ws.Cells.Item(144, 14).Value = 3.75
$ws.Cells.Item(144, 15).Value = 2.1
$ws.Cells.Item(144, 16).Value = 0.25
$ws.Cells.Item(144, 17).Value = 1.875
$ws.Cells.Item(144, 18).Value = 1.925
$ws.Cells.Item(144, 19).Value = 3.5
$ws.Cells.Item(144, 20).Value = 1.875
$ws.Cells.Item(144, 21).Value = 1.925
$ws.Cells.Item(144, 22).Value = 1.75
$ws.Cells.Item(144, 23).Value = -1
$ws.Cells.Item(144, 25).Value = 0.875
$ws.Cells.Item(144, 26).Value = -1
$ws.Cells.Item(144, 28).Value = 0.925
$ws.Cells.Item(151, 2).Value = 8161940
$ws.Cells.Item(151, 5).Value = 'Germania Schneiche'
$ws.Cells.Item(151, 6).Value = 'TuS Sachsenhausen'
$ws.Cells.Item(151, 7).Value = 1
$ws.Cells.Item(151, 8).Value = 4
$ws.Cells.Item(151, 10).Value = 3.25
$ws.Cells.Item(151, 12).Value = 1.8
$ws.Cells.Item(151, 13).Value = 3.25
$ws.Cells.Item(151, 14).Value = 4
$ws.Cells.Item(151, 15).Value = 1.8
$ws.Cells.Item(151, 16).Value = 0.5
$ws.Cells.Item(151, 17).Value = 1.95
$ws.Cells.Item(151, 18).Value = 1.85
$ws.Cells.Item(151, 19).Value = 3.5
$ws.Cells.Item(151, 20).Value = 1.9
$ws.Cells.Item(151, 21).Value = 1.9
$ws.Cells.Item(151, 24).Value = 0.8
$ws.Cells.Item(151, 26).Value = 0.8500000000000001
$ws.Cells.Item(151, 27).Value = 0.8999999999999999
$ws.Cells.Item(151, 28).Value = -1
$ws.Cells.Item(152, 2).Value = 8162017
$ws.Cells.Item(152, 5).Value = 'VfR Fehlheim'
$ws.Cells.Item(152, 6).Value = 'SV Pars NeuIsenburg'
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 3
$ws.Cells.Item(152, 10).Value = 3.5
$ws.Cells.Item(152, 12).Value = 1.727
$ws.Cells.Item(152, 13).Value = 5
$ws.Cells.Item(152, 14).Value = 4.75
$ws.Cells.Item(152, 15).Value = 1.4
$ws.Cells.Item(152, 16).Value = 1.25
$ws.Cells.Item(152, 17).Value = 1.9
$ws.Cells.Item(152, 18).Value = 1.9
$ws.Cells.Item(152, 19).Value = 4
$ws.Cells.Item(152, 20).Value = 1.925
$ws.Cells.Item(152, 21).Value = 1.875
$ws.Cells.Item(152, 24).Value = 0.3999999999999999
$ws.Cells.Item(152, 26).Value = 0.8999999999999999
$ws.Cells.Item(152, 27).Value = -1
$ws.Cells.Item(152, 28).Value = 0.875

# --- Append new row 156 (new match result) ---
$ws.Cells.Item(156, 1).Value = 154
$ws.Cells.Item(156, 2).Value = 8191504
$ws.Cells.Item(156, 3).Value = 'Germany Verbandsliga'
$ws.Cells.Item(156, 4).Value = 45420.59375
$ws.Cells.Item(156, 5).Value = '1 FC Frankfurt'
$ws.Cells.Item(156, 6).Value = 'FV Preussen Eberswalde'
$ws.Cells.Item(156, 7).Value = 2
$ws.Cells.Item(156, 8).Value = 2
$ws.Cells.Item(156, 9).Value = 'D'
$ws.Cells.Item(156, 10).Value = 1.062
$ws.Cells.Item(156, 11).Value = 9.5
$ws.Cells.Item(156, 12).Value = 15
$ws.Cells.Item(156, 13).Value = 1.111
$ws.Cells.Item(156, 14).Value = 8.5
$ws.Cells.Item(156, 15).Value = 17
$ws.Cells.Item(156, 16).Value = -2.75
$ws.Cells.Item(156, 17).Value = 1.85
$ws.Cells.Item(156, 18).Value = 1.95
$ws.Cells.Item(156, 19).Value = 4
$ws.Cells.Item(156, 20).Value = 1.9
$ws.Cells.Item(156, 21).Value = 1.9
$ws.Cells.Item(156, 22).Value = -1
$ws.Cells.Item(156, 23).Value = 7.5
$ws.Cells.Item(156, 24).Value = -1
$ws.Cells.Item(156, 25).Value = -1
$ws.Cells.Item(156, 26).Value = 0.95
$ws.Cells.Item(156, 27).Value = 0
$ws.Cells.Item(156, 28).Value = 0

# Copy cell formatting (styles) from the last existing data row (155) to the new row 156
$ws.Cells.Item(155, 1).Copy() | Out-Null
$ws.Cells.Item(156, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(155, 4).Copy() | Out-Null
$ws.Cells.Item(156, 4).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
